$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33:F33").Copy()
$ws.Range("A34:F34").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
